# Insert a new data row at row 124 (pushing the existing rows 124-239 down to
# 125-240, extending the table from A1:R239 to A1:R240), then populate the new
# row with its own data (a newer weekly price observation for Ciboulette /
# Femacal de La Calera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a full row above the current row 124; this shifts every row from
# 124..239 down by one (to 125..240) and carries their formatting along.
$ws.Rows("124:124").Insert()

# Fill in the freshly inserted row 124 with the new observation.
$ws.Range("A124").Value = 3
$ws.Range("B124").Value = "Femacal de La Calera"
$ws.Range("C124").Value = "Coquimbo"
$ws.Range("D124").Value = 44566
$ws.Range("E124").Value = 5
$ws.Range("F124").Value = 100112039
$ws.Range("G124").Value = "Ciboulette"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 130
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 1500
$ws.Range("M124").Value = 1500
$ws.Range("N124").Value = "$/docena de atados"
$ws.Range("O124").Value = "Provincia de Quillota"
$ws.Range("P124").Value = 500
$ws.Range("Q124").Value = 3
$ws.Range("R124").Value = "Hortaliza"
